$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the spacer/press-in-thread line item (row 5, "9774040151R") entirely.
# Deleting the whole row shifts the rows below it (PicoBlade connector, RGB LED,
# keyswitch rows) up by one, matching the new BOM layout.
$ws.Rows("5").Delete()
